$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets("ALC")
$ws.Range("H19").Value = 1170.8334
$ws.Range("I19").Value = 938.4286
$ws.Range("K19").Value = 938.4286
$ws.Range("M19").Value = -763.4286

$ws.Range("H28").Value = 3072.647
$ws.Range("I28").Value = 1951.1428
$ws.Range("J28").Value = 8306.333000000001
$ws.Range("K28").Value = 1951.1428
$ws.Range("L28").Value = 8306.333000000001
$ws.Range("M28").Value = -1466.1428
$ws.Range("N28").Value = -9276.333000000001

$ws.Range("H32").Value = 102
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 102
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 102
$ws.Range("N32").Value = -754
$ws.Range("M32").ClearContents()

$ws.Range("H70").Value = 90278750
$ws.Range("J70").Value = 92593610
$ws.Range("L70").Value = 277780830
$ws.Range("N70").Value = -277781370

$ws.Range("H73").Value = 90278750
$ws.Range("J73").Value = 92593610
$ws.Range("L73").Value = 277780830
$ws.Range("N73").Value = -277782702

$ws.Range("H76").Value = 10095.75
$ws.Range("I76").Value = 9984.727999999999
$ws.Range("K76").Value = 9984.727999999999
$ws.Range("M76").Value = -9669.727999999999

$ws.Range("H79").Value = 10095.75
$ws.Range("I79").Value = 9984.727999999999
$ws.Range("K79").Value = 9984.727999999999
$ws.Range("M79").Value = -8892.727999999999

$ws.Range("H98").Value = 2492.7173
$ws.Range("I98").Value = 2528.1538
$ws.Range("J98").Value = 2295.2856
$ws.Range("K98").Value = 2528.1538
$ws.Range("L98").Value = 2295.2856
$ws.Range("M98").Value = -1030.1538
$ws.Range("N98").Value = -5291.2856

$ws.Range("H106").Value = 2261.4
$ws.Range("I106").Value = 2261.4
$ws.Range("K106").Value = 2261.4
$ws.Range("M106").Value = -1630.4

$ws.Range("H113").Value = 38202820
$ws.Range("I113").Value = 15879302
$ws.Range("K113").Value = 15879302
$ws.Range("M113").Value = -15876048

$ws.Range("H116").Value = 125010000
$ws.Range("I116").Value = 250000000
$ws.Range("J116").Value = 20000
$ws.Range("K116").Value = 250000000
$ws.Range("L116").Value = 20000
$ws.Range("M116").Value = -249996558
$ws.Range("N116").Value = -26884

$ws.Range("H118").Value = 542
$ws.Range("I118").Value = 542
$ws.Range("J118").Value = 0
$ws.Range("K118").Value = 1626
$ws.Range("L118").Value = 0
$ws.Range("M118").Value = 31
$ws.Range("N118").ClearContents()

$ws.Range("H122").Value = 2492.7173
$ws.Range("I122").Value = 2528.1538
$ws.Range("J122").Value = 2295.2856
$ws.Range("K122").Value = 7584.4614
$ws.Range("L122").Value = 6885.8568
$ws.Range("M122").Value = -5134.4614
$ws.Range("N122").Value = -11785.8568

$ws.Range("H132").Value = 1521.3
$ws.Range("I132").Value = 1435.8276
$ws.Range("K132").Value = 4307.4828
$ws.Range("M132").Value = -1777.4828

$ws.Range("H138").Value = 2131461.2
$ws.Range("J138").Value = 3453023
$ws.Range("L138").Value = 10359069
$ws.Range("N138").Value = -10369349

$ws = $wb.Worksheets("ARM")
$ws.Range("H32").Value = 2454785.5
$ws.Range("I32").Value = 2553954
$ws.Range("K32").Value = 2553954
$ws.Range("M32").Value = -2553667

$ws.Range("H61").Value = 8931.629999999999
$ws.Range("I61").Value = 3921.1428
$ws.Range("J61").Value = 10685.3
$ws.Range("K61").Value = 3921.1428
$ws.Range("L61").Value = 10685.3
$ws.Range("M61").Value = -3709.1428
$ws.Range("N61").Value = -11109.3

$ws.Range("H74").Value = 33352.793
$ws.Range("I74").Value = 44401.332
$ws.Range("K74").Value = 44401.332
$ws.Range("M74").Value = -43527.332

$ws.Range("H77").Value = 33352.793
$ws.Range("I77").Value = 44401.332
$ws.Range("K77").Value = 222006.66
$ws.Range("M77").Value = -217638.66

$ws.Range("H122").Value = 3218.4546
$ws.Range("I122").Value = 2891.7273
$ws.Range("J122").Value = 3545.182
$ws.Range("K122").Value = 8675.1819
$ws.Range("L122").Value = 10635.546
$ws.Range("M122").Value = -6225.1819
$ws.Range("N122").Value = -15535.546

$ws.Range("H136").Value = 8931.629999999999
$ws.Range("I136").Value = 3921.1428
$ws.Range("J136").Value = 10685.3
$ws.Range("K136").Value = 11763.4284
$ws.Range("L136").Value = 32055.9
$ws.Range("M136").Value = -9213.428400000001
$ws.Range("N136").Value = -37155.89999999999

$ws = $wb.Worksheets("BSM")
$ws.Range("H113").Value = 7513.5
$ws.Range("I113").Value = 7513.5
$ws.Range("K113").Value = 7513.5
$ws.Range("M113").Value = -5343.5

$ws.Range("H134").Value = 9524.1875
$ws.Range("I134").Value = 4056.1428
$ws.Range("K134").Value = 12168.4284
$ws.Range("M134").Value = -9633.428400000001

$ws = $wb.Worksheets("CRP")
$ws.Range("H31").Value = 12835.353
$ws.Range("I31").Value = 4201.1665
$ws.Range("K31").Value = 4201.1665
$ws.Range("M31").Value = -3906.1665

$ws.Range("H34").Value = 12835.353
$ws.Range("I34").Value = 4201.1665
$ws.Range("K34").Value = 4201.1665
$ws.Range("M34").Value = -3999.1665

$ws.Range("H105").Value = 7940357
$ws.Range("I105").Value = 11905545
$ws.Range("K105").Value = 11905545
$ws.Range("M105").Value = -11903798

$ws.Range("H122").Value = 2859
$ws.Range("I122").Value = 2618.625
$ws.Range("K122").Value = 7855.875
$ws.Range("M122").Value = -5405.875

$ws = $wb.Worksheets("CUL")
$ws.Range("H56").Value = 7247.5
$ws.Range("I56").Value = 7247.5
$ws.Range("K56").Value = 7247.5
$ws.Range("M56").Value = -6717.5

$ws.Range("H68").Value = 2281.48
$ws.Range("J68").Value = 2511.2778
$ws.Range("L68").Value = 7533.8334
$ws.Range("N68").Value = -9155.8334

$ws.Range("H71").Value = 2281.48
$ws.Range("J71").Value = 2511.2778
$ws.Range("L71").Value = 22601.5002
$ws.Range("N71").Value = -30713.5002

$ws.Range("H98").Value = 3996.25
$ws.Range("I98").Value = 2995
$ws.Range("J98").Value = 4997.5
$ws.Range("K98").Value = 8985
$ws.Range("L98").Value = 14992.5
$ws.Range("M98").Value = -7487
$ws.Range("N98").Value = -17988.5

$ws.Range("H132").Value = 6404.1763
$ws.Range("J132").Value = 8648.375
$ws.Range("L132").Value = 77835.375
$ws.Range("N132").Value = -82895.375

$ws = $wb.Worksheets("LTW")
$ws.Range("H7").Value = 7195.7393
$ws.Range("I7").Value = 4374
$ws.Range("K7").Value = 4374
$ws.Range("M7").Value = -4262

$ws.Range("H98").Value = 54804
$ws.Range("J98").Value = 54804
$ws.Range("L98").Value = 54804
$ws.Range("N98").Value = -60794

$ws.Range("H122").Value = 9455.5
$ws.Range("I122").Value = 12414.5
$ws.Range("J122").Value = 7680.1
$ws.Range("K122").Value = 37243.5
$ws.Range("L122").Value = 23040.3
$ws.Range("M122").Value = -34793.5
$ws.Range("N122").Value = -27940.3

$ws.Range("H126").Value = 7195.7393
$ws.Range("I126").Value = 4374
$ws.Range("K126").Value = 13122
$ws.Range("M126").Value = -10652

$ws = $wb.Worksheets("WVR")
$ws.Range("H81").Value = 12357723
$ws.Range("I81").Value = 715699
$ws.Range("J81").Value = 66687170
$ws.Range("K81").Value = 1431398
$ws.Range("L81").Value = 133374340
$ws.Range("M81").Value = -1430337
$ws.Range("N81").Value = -133376462

$ws.Range("H84").Value = 12357723
$ws.Range("I84").Value = 715699
$ws.Range("J84").Value = 66687170
$ws.Range("K84").Value = 7156990
$ws.Range("L84").Value = 666871700
$ws.Range("M84").Value = -7151686
$ws.Range("N84").Value = -666882308
